# Refresh the crypto Price (D) and Volume(1h) (E) columns with the latest values
# pulled from coinranking.com, as produced by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.727.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.70%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.401.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.02%  "

$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.22%  "

$ws.Range("E8").Value = "  -1.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.381.15"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("E11").Value = "  +0.34%  "

$ws.Range("E12").Value = "  -3.36%  "

$ws.Range("E13").Value = "  -2.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.57%  "

$ws.Range("E15").Value = "  -2.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.826.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.635.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.377.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.88%  "

$ws.Range("E22").Value = "  -2.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.06%  "

$ws.Range("E24").Value = "  +0.26%  "

$ws.Range("E25").Value = "  -6.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "576.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = ("0.0" + ([string][char]0x2083) + "0904")
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.50%  "

$ws.Range("E32").Value = "  -7.18%  "

$ws.Range("E34").Value = "  -8.01%  "

$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("E36").Value = "  -6.88%  "

$ws.Range("E37").Value = "  -3.50%  "

$ws.Range("E38").Value = "  -4.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "146.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.96%  "

$ws.Range("E40").Value = "  -1.48%  "

$ws.Range("E41").Value = "  -5.04%  "

$ws.Range("E42").Value = "  +0.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.26%  "

$ws.Range("E44").Value = "  -5.50%  "

$ws.Range("E45").Value = "  -5.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = ("0.0" + ([string][char]0x2086) + "0285")
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +20.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.08%  "

$ws.Range("E48").Value = "  -4.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.583"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.23%  "

$ws.Range("E50").Value = "  -4.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.92%  "
